# Add optimization log flag and conditional loading of scale values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: climate_change_factor_gnrl_hydropower_availability -> updated scale value
$ws.Range("J4:AS4").Value = 0.6193541768728235

# Row 5: elasticity_gnrl_rate_occupancy_to_gdppc -> updated scale value
$ws.Range("J5:AS5").Value = -0.07866934995703223
